# Removed reference to Dev Hub
# The "Attributes" sheet (sheet2) has a row (row 57) whose value is the
# "app" = "DevHub" entry for user "j_thomas". That row is removed entirely,
# causing rows below it to shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# Select the row to make the resulting selection match what Excel leaves
# behind after deleting a row (entire row selection).
$ws.Rows.Item(57).Select()

# Delete the entire row 57 (the DevHub entry), shifting rows 58:67 up.
$ws.Rows.Item(57).Delete()
